# Generate Report for Handoff
# - Flip the localization "Status" from "In Translation" to "Ready for handoff"
#   on every sheet that tracks it (Overview's per-language columns + each
#   language sheet's Status column), and bump the associated timestamps to
#   reflect the new handoff-generation run.
# - The Status text grew longer, so the Status column(s) get re-autofit to a
#   wider width on every sheet that shows them.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ------------------
$overview.Range("E2").Value = "Ready for handoff"   # zh-cn status column
$overview.Range("F2").Value = "Ready for handoff"   # de-de status column
$zhcn.Range("C2").Value     = "Ready for handoff"   # Status column
$dede.Range("C2").Value     = "Ready for handoff"   # Status column

# --- Timestamps refreshed by the new handoff-report generation run ---------
$overview.Range("G2").Value = "2016-08-12 07:12:33"   # Latest HO Xliff Generate Date
$zhcn.Range("H2").Value     = "2016-08-12 07:12:26"   # Latest Handoff Datetime (zh-cn)
$dede.Range("H2").Value     = "2016-08-12 07:12:33"   # Latest Handoff Datetime (de-de)

# --- Re-autofit the Status column(s) now that the text is wider ------------
# (ColumnWidth is rounded to the host's pixel grid, so we feed it the input
# that lands on the grid cell closest to the generator's computed width.)
$overview.Columns.Item(5).ColumnWidth = 16.25   # E: zh-cn Status
$overview.Columns.Item(6).ColumnWidth = 16.25   # F: de-de Status
$zhcn.Columns.Item(3).ColumnWidth     = 16.25   # C: Status
$dede.Columns.Item(3).ColumnWidth     = 16.25   # C: Status
